# Auto-generated edit script: update crypto price/volume table
# Matches the diff: updates D (Price) and E (Volume 1h) columns for rows 2-51,
# and swaps the OKB / Fetch.AI rows (39 <-> 40).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.301.90"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.91%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'2.601.44"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.06%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'510.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.15%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'154.53"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.42%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.587"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.09%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'2.614.10"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.38%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'6.66"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.49%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("E11").Value = "  -0.51%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.345"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.67%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("E13").Value = "  +1.66%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'3.056.60"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.94%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'60.314.54"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.97%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("E16").Value = "  -1.07%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("E17").Value = "  +0.24%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'2.602.08"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.10%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'4.75"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.12%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'350.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.74%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("E21").Value = "  +0.45%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'6.12"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.14%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'60.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("E25").Value = "  -0.68%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'0.166"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.44%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D28").Value = "0.0₃0837"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.58%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("E29").Value = "  -2.27%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("E30").Value = "  -0.14%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'19.38"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.24%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'151.22"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.13%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("E33").Value = "  -0.82%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("E34").Value = "  +0.27%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'3.99"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.23%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("E36").Value = "  -2.81%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'0.876"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.64%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("E38").Value = "  -1.94%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("E41").Value = "  -0.44%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'293.41"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.86%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("E43").Value = "  -3.72%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("E44").Value = "  -0.25%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.998"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.19%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'0.0552"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.92%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'19.66"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.37%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'4.87"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.74%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'0.0233"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.32%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("E50").Value = "  +0.00%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'1.996.02"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.38%  "
$ws.Range("E51").Style = "Normal"

# Rows 39/40: OKB and Fetch.AI swapped places in the ranking
$ws.Range("B39").Value = "OKB"
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = "'36.24"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.29%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("B40").Value = "Fetch.AI"
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = "'0.841"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.44%  "
$ws.Range("E40").Style = "Normal"
